# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ROAN")

# Row 21 - Earnings Before Interest And Taxes: replace "NA" placeholders with actual figures
$ws.Range("D21").Value = 3319900
$ws.Range("E21").Value = 29700
$ws.Range("F21").Value = -3196600
$ws.Range("G21").Value = 626500
$ws.Range("H21").Value = 135800
$ws.Range("I21").Value = 222300

# Row 41 - Cash And Cash Equivalents
$ws.Range("D41").Value = 1500

# Row 43 - Net Receivables
$ws.Range("D43").Value = 220000

# Row 45 - Other Current Assets
$ws.Range("D45").Value = 246200

# Row 46 - Total Current Assets
$ws.Range("D46").Value = 82100

# Row 48 - Property Plant and Equipment
$ws.Range("D48").Value = 2251900

# Row 52 - Other Assets
$ws.Range("D52").Value = 209600

# Row 54 - Total Assets
$ws.Range("D54").Value = 1885600

# Row 57 - Accounts Payable
$ws.Range("D57").Value = 194100

# Row 59 - Other Current Liabilities
$ws.Range("D59").Value = 111200

# Row 60 - Total Current Liabilities
$ws.Range("D60").Value = 203300

# Row 61 - Long Term Debt
$ws.Range("D61").Value = 85300

# Row 62 - Other Liabilities
$ws.Range("D62").Value = 12100

# Row 66 - Total Liabilities
$ws.Range("D66").Value = 300800

# Row 72 - Retained Earnings
$ws.Range("D72").Value = 0

# Row 76 - Total Stockholder Equity
$ws.Range("D76").Value = 1584800

# Row 101 - Other Cash Flows from Financing Activities: J column becomes NA
$ws.Range("J101").Value = "NA"
